# Updated cryptos list: refresh the Price (D) and Volume(1h) (E) columns
# for rows 2-51 with the latest scraped figures. Row 28 is unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some new Price values (single-dot decimals such as "594.31") would
# otherwise be auto-converted to numbers by Excel's type inference when
# assigned through .Value. Pre-format those specific cells as Text so the
# literal string is preserved, matching the source data (which keeps
# thousands-dot-formatted prices, e.g. "66.877.05", as plain text too).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

$ws.Range('D2').Value = '66.877.05'
$ws.Range('E2').Value = '  +0.49%  '
$ws.Range('D3').Value = '3.498.26'
$ws.Range('E3').Value = '  +0.30%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '594.31'
$ws.Range('E5').Value = '  +0.37%  '
$ws.Range('D6').Value = '172.50'
$ws.Range('E6').Value = '  +2.22%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  -0.35%  '
$ws.Range('E9').Value = '  +3.54%  '
$ws.Range('D10').Value = '7.19'
$ws.Range('E10').Value = '  -1.79%  '
$ws.Range('D11').Value = '0.431'
$ws.Range('E11').Value = '  -1.20%  '
$ws.Range('D12').Value = '4.102.17'
$ws.Range('E12').Value = '  +0.37%  '
$ws.Range('E13').Value = '  -0.33%  '
$ws.Range('D14').Value = '29.23'
$ws.Range('E14').Value = '  +3.83%  '
$ws.Range('D15').Value = '66.883.35'
$ws.Range('E15').Value = '  +0.45%  '
$ws.Range('E16').Value = '  +0.55%  '
$ws.Range('D17').Value = '3.495.37'
$ws.Range('E17').Value = '  +0.23%  '
$ws.Range('E18').Value = '  -0.32%  '
$ws.Range('D19').Value = '14.26'
$ws.Range('E19').Value = '  +1.82%  '
$ws.Range('D20').Value = '395.19'
$ws.Range('E20').Value = '  +0.95%  '
$ws.Range('D21').Value = '7.95'
$ws.Range('E21').Value = '  +0.63%  '
$ws.Range('D22').Value = '73.43'
$ws.Range('E22').Value = '  +0.71%  '
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').Value = '0.535'
$ws.Range('E24').Value = '  +0.23%  '
$ws.Range('E25').Value = '  -0.29%  '
$ws.Range('D26').Value = '10.22'
$ws.Range('E26').Value = '  -0.67%  '
$ws.Range('E27').Value = '  +0.40%  '
$ws.Range('D29').Value = '6.16'
$ws.Range('E29').Value = '  -2.27%  '
$ws.Range('E30').Value = '  -2.03%  '
$ws.Range('E31').Value = '  -0.24%  '
$ws.Range('E32').Value = '  +0.50%  '
$ws.Range('D33').Value = '7.36'
$ws.Range('E33').Value = '  -0.36%  '
$ws.Range('E34').Value = '  +0.37%  '
$ws.Range('D35').Value = '162.63'
$ws.Range('E35').Value = '  +0.66%  '
$ws.Range('E36').Value = '  -1.07%  '
$ws.Range('E37').Value = '  -0.60%  '
$ws.Range('E38').Value = '  +1.77%  '
$ws.Range('E39').Value = '  +0.45%  '
$ws.Range('E40').Value = '  -0.71%  '
$ws.Range('D41').Value = '2.836.03'
$ws.Range('E41').Value = '  +2.46%  '
$ws.Range('D42').Value = '27.10'
$ws.Range('E42').Value = '  +1.06%  '
$ws.Range('D43').Value = '26.17'
$ws.Range('E43').Value = '  -1.09%  '
$ws.Range('D44').Value = '42.71'
$ws.Range('E44').Value = '  -0.99%  '
$ws.Range('E45').Value = '  +2.35%  '
$ws.Range('D46').Value = '0.0303'
$ws.Range('E46').Value = '  -3.02%  '
$ws.Range('D47').Value = '336.67'
$ws.Range('E47').Value = '  -2.58%  '
$ws.Range('D48').Value = '34.65'
$ws.Range('E48').Value = '  +2.26%  '
$ws.Range('E49').Value = '  -1.12%  '
$ws.Range('E50').Value = '  -1.10%  '
$ws.Range('E51').Value = '  -3.92%  '
